# Scheduled runner update: refresh market-board derived price/profit figures
# across the leve-profit tables on each crafting-class sheet (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Only the price/profit columns (H-N) for the
# affected rows are touched; all other data is left as-is.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 10000
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -10586
$ws.Range("H17").Value = 1856921.2
$ws.Range("J17").Value = 1856921.2
$ws.Range("L17").Value = 5570763.6
$ws.Range("N17").Value = -5571099.6
$ws.Range("H112").Value = 1148.0465
$ws.Range("J112").Value = 1148.0465
$ws.Range("L112").Value = 3444.1395
$ws.Range("N112").Value = -5660.139499999999
$ws.Range("H126").Value = 11996.25
$ws.Range("J126").Value = 11996.25
$ws.Range("L126").Value = 11996.25
$ws.Range("N126").Value = -21876.25
$ws.Range("H138").Value = 3153.68
$ws.Range("I138").Value = 1515.34
$ws.Range("J138").Value = 4792.02
$ws.Range("K138").Value = 4546.02
$ws.Range("L138").Value = 14376.06
$ws.Range("M138").Value = 593.9800000000005
$ws.Range("N138").Value = -24656.06

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2650
$ws.Range("I3").Value = 300
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 300
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -185
$ws.Range("N3").Value = -5230
$ws.Range("H102").Value = 2099.6072
$ws.Range("I102").Value = 2099.6072
$ws.Range("K102").Value = 2099.6072
$ws.Range("M102").Value = -477.6071999999999
$ws.Range("H134").Value = 40505.8
$ws.Range("J134").Value = 45632.25
$ws.Range("L134").Value = 45632.25
$ws.Range("N134").Value = -55772.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2627.0588
$ws.Range("I105").Value = 2562.7856
$ws.Range("J105").Value = 2927
$ws.Range("K105").Value = 2562.7856
$ws.Range("L105").Value = 2927
$ws.Range("M105").Value = -815.7856000000002
$ws.Range("N105").Value = -6421
$ws.Range("H107").Value = 544.381
$ws.Range("I107").Value = 409.5
$ws.Range("J107").Value = 976
$ws.Range("K107").Value = 409.5
$ws.Range("L107").Value = 976
$ws.Range("M107").Value = 1510.5
$ws.Range("N107").Value = -4816
$ws.Range("H140").Value = 48864.168
$ws.Range("J140").Value = 48864.168
$ws.Range("L140").Value = 48864.168
$ws.Range("N140").Value = -59224.168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 10000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10278
$ws.Range("H31").Value = 3205.087
$ws.Range("I31").Value = 1370.9333
$ws.Range("J31").Value = 6644.125
$ws.Range("K31").Value = 1370.9333
$ws.Range("L31").Value = 6644.125
$ws.Range("M31").Value = -1075.9333
$ws.Range("N31").Value = -7234.125
$ws.Range("H34").Value = 3205.087
$ws.Range("I34").Value = 1370.9333
$ws.Range("J34").Value = 6644.125
$ws.Range("K34").Value = 1370.9333
$ws.Range("L34").Value = 6644.125
$ws.Range("M34").Value = -1168.9333
$ws.Range("N34").Value = -7048.125
$ws.Range("H107").Value = 418.8125
$ws.Range("I107").Value = 346.1
$ws.Range("J107").Value = 540
$ws.Range("K107").Value = 346.1
$ws.Range("L107").Value = 540
$ws.Range("M107").Value = 1573.9
$ws.Range("N107").Value = -4380
$ws.Range("H138").Value = 45000
$ws.Range("J138").Value = 45000
$ws.Range("L138").Value = 45000
$ws.Range("N138").Value = -55280
$ws.Range("H140").Value = 103226
$ws.Range("J140").Value = 103226
$ws.Range("L140").Value = 103226
$ws.Range("N140").Value = -113586

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 617.8
$ws.Range("I10").Value = 311.14285
$ws.Range("J10").Value = 1333.3334
$ws.Range("K10").Value = 933.4285500000001
$ws.Range("L10").Value = 4000.0002
$ws.Range("M10").Value = -794.4285500000001
$ws.Range("N10").Value = -4278.0002
$ws.Range("H122").Value = 1268.4138
$ws.Range("I122").Value = 504
$ws.Range("J122").Value = 3670.8572
$ws.Range("K122").Value = 4536
$ws.Range("L122").Value = 33037.7148
$ws.Range("M122").Value = -2086
$ws.Range("N122").Value = -37937.7148
$ws.Range("H131").Value = 808.7553
$ws.Range("I131").Value = 327.375
$ws.Range("J131").Value = 853.5349
$ws.Range("K131").Value = 982.125
$ws.Range("L131").Value = 2560.6047
$ws.Range("M131").Value = 4057.875
$ws.Range("N131").Value = -12640.6047

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 872.0833
$ws.Range("I2").Value = 953.1818
$ws.Range("J2").Value = 803.46155
$ws.Range("K2").Value = 953.1818
$ws.Range("L2").Value = 803.46155
$ws.Range("M2").Value = -840.1818
$ws.Range("N2").Value = -1029.46155
$ws.Range("H9").Value = 500
$ws.Range("I9").Value = 500
$ws.Range("K9").Value = 500
$ws.Range("M9").Value = -330
$ws.Range("H109").Value = 29000
$ws.Range("J109").Value = 29000
$ws.Range("L109").Value = 29000
$ws.Range("N109").Value = -31080
$ws.Range("H135").Value = 38632.223
$ws.Range("J135").Value = 38632.223
$ws.Range("L135").Value = 38632.223
$ws.Range("N135").Value = -48772.223
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2600.72
$ws.Range("I100").Value = 2553.8262
$ws.Range("J100").Value = 3140
$ws.Range("K100").Value = 2553.8262
$ws.Range("L100").Value = 3140
$ws.Range("M100").Value = -2012.8262
$ws.Range("N100").Value = -4222
$ws.Range("H127").Value = 42318.332
$ws.Range("J127").Value = 42318.332
$ws.Range("L127").Value = 42318.332
$ws.Range("N127").Value = -52238.332
$ws.Range("H136").Value = 2730.7407
$ws.Range("I136").Value = 2445.3044
$ws.Range("J136").Value = 4372
$ws.Range("K136").Value = 7335.9132
$ws.Range("L136").Value = 13116
$ws.Range("M136").Value = -4785.9132
$ws.Range("N136").Value = -18216

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 49777
$ws.Range("J128").Value = 49777
$ws.Range("L128").Value = 49777
$ws.Range("N128").Value = -59737
$ws.Range("H136").Value = 3523.255
$ws.Range("I136").Value = 543.9375
$ws.Range("J136").Value = 8541.053
$ws.Range("K136").Value = 1631.8125
$ws.Range("L136").Value = 25623.159
$ws.Range("M136").Value = 918.1875
$ws.Range("N136").Value = -30723.159
